$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell, new literal text, and whether the text must be
# protected from Excel's automatic number/date inference (values such as
# "1.001" or "0.2870" would otherwise be re-typed as numeric and lose their
# exact textual form).
$updates = @(
    @{ Cell = 'D2'; Value = '24.951.94'; ForceText = $false },
    @{ Cell = 'E2'; Value = '  +2.13%  '; ForceText = $false },
    @{ Cell = 'D3'; Value = '1.699.08'; ForceText = $false },
    @{ Cell = 'E3'; Value = '  +0.75%  '; ForceText = $false },
    @{ Cell = 'D4'; Value = '1.001'; ForceText = $true },
    @{ Cell = 'E4'; Value = '  +0.13%  '; ForceText = $false },
    @{ Cell = 'D5'; Value = '315.03'; ForceText = $true },
    @{ Cell = 'E5'; Value = '  -0.26%  '; ForceText = $false },
    @{ Cell = 'E6'; Value = '  +0.18%  '; ForceText = $false },
    @{ Cell = 'D7'; Value = '0.3973'; ForceText = $true },
    @{ Cell = 'E7'; Value = '  +1.69%  '; ForceText = $false },
    @{ Cell = 'D8'; Value = '0.4021'; ForceText = $true },
    @{ Cell = 'E8'; Value = '  -0.13%  '; ForceText = $false },
    @{ Cell = 'D9'; Value = '1.463'; ForceText = $true },
    @{ Cell = 'E9'; Value = '  -1.70%  '; ForceText = $false },
    @{ Cell = 'D10'; Value = '52.78'; ForceText = $true },
    @{ Cell = 'E10'; Value = '  +0.57%  '; ForceText = $false },
    @{ Cell = 'D11'; Value = '1.003'; ForceText = $true },
    @{ Cell = 'E11'; Value = '  +0.31%  '; ForceText = $false },
    @{ Cell = 'D12'; Value = '0.08795'; ForceText = $true },
    @{ Cell = 'E12'; Value = '  +0.21%  '; ForceText = $false },
    @{ Cell = 'D13'; Value = '26.06'; ForceText = $true },
    @{ Cell = 'E13'; Value = '  -2.12%  '; ForceText = $false },
    @{ Cell = 'D14'; Value = '7.454'; ForceText = $true },
    @{ Cell = 'E14'; Value = '  -0.04%  '; ForceText = $false },
    @{ Cell = 'D15'; Value = '0.00001351'; ForceText = $true },
    @{ Cell = 'E15'; Value = '  +0.39%  '; ForceText = $false },
    @{ Cell = 'D16'; Value = '7.943'; ForceText = $true },
    @{ Cell = 'E16'; Value = '  -2.41%  '; ForceText = $false },
    @{ Cell = 'D17'; Value = '1.707.20'; ForceText = $false },
    @{ Cell = 'E17'; Value = '  +1.43%  '; ForceText = $false },
    @{ Cell = 'D18'; Value = '95.94'; ForceText = $true },
    @{ Cell = 'E18'; Value = '  -2.09%  '; ForceText = $false },
    @{ Cell = 'D19'; Value = '0.07184'; ForceText = $true },
    @{ Cell = 'E19'; Value = '  -0.71%  '; ForceText = $false },
    @{ Cell = 'D20'; Value = '20.58'; ForceText = $true },
    @{ Cell = 'E20'; Value = '  +1.41%  '; ForceText = $false },
    @{ Cell = 'D21'; Value = '7.318'; ForceText = $true },
    @{ Cell = 'E21'; Value = '  +0.64%  '; ForceText = $false },
    @{ Cell = 'D22'; Value = '1.001'; ForceText = $true },
    @{ Cell = 'E22'; Value = '  +0.18%  '; ForceText = $false },
    @{ Cell = 'D23'; Value = '14.36'; ForceText = $true },
    @{ Cell = 'E23'; Value = '  +0.70%  '; ForceText = $false },
    @{ Cell = 'D24'; Value = '24.942.51'; ForceText = $false },
    @{ Cell = 'E24'; Value = '  +2.13%  '; ForceText = $false },
    @{ Cell = 'D25'; Value = '2.358'; ForceText = $true },
    @{ Cell = 'E25'; Value = '  +0.95%  '; ForceText = $false },
    @{ Cell = 'D26'; Value = '2.938'; ForceText = $true },
    @{ Cell = 'E26'; Value = '  -3.34%  '; ForceText = $false },
    @{ Cell = 'D27'; Value = '23.75'; ForceText = $true },
    @{ Cell = 'E27'; Value = '  +5.18%  '; ForceText = $false },
    @{ Cell = 'D28'; Value = '6.168'; ForceText = $true },
    @{ Cell = 'E28'; Value = '  +15.61%  '; ForceText = $false },
    @{ Cell = 'D29'; Value = '161.89'; ForceText = $true },
    @{ Cell = 'E29'; Value = '  -3.15%  '; ForceText = $false },
    @{ Cell = 'D30'; Value = '150.91'; ForceText = $true },
    @{ Cell = 'E30'; Value = '  +9.30%  '; ForceText = $false },
    @{ Cell = 'D31'; Value = '8.343'; ForceText = $true },
    @{ Cell = 'E31'; Value = '  -1.67%  '; ForceText = $false },
    @{ Cell = 'D32'; Value = '2.646'; ForceText = $true },
    @{ Cell = 'E32'; Value = '  +26.40%  '; ForceText = $false },
    @{ Cell = 'D33'; Value = '1.895.81'; ForceText = $false },
    @{ Cell = 'E33'; Value = '  +1.43%  '; ForceText = $false },
    @{ Cell = 'D34'; Value = '0.08553'; ForceText = $true },
    @{ Cell = 'E34'; Value = '  -2.49%  '; ForceText = $false },
    @{ Cell = 'D35'; Value = '0.03147'; ForceText = $true },
    @{ Cell = 'E35'; Value = '  +4.18%  '; ForceText = $false },
    @{ Cell = 'B36'; Value = 'InternetComputer(DFINITY)'; ForceText = $false },
    @{ Cell = 'C36'; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; ForceText = $false },
    @{ Cell = 'D36'; Value = '7.152'; ForceText = $true },
    @{ Cell = 'E36'; Value = '  -1.76%  '; ForceText = $false },
    @{ Cell = 'B37'; Value = 'ImmutableX'; ForceText = $false },
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; ForceText = $false },
    @{ Cell = 'D37'; Value = '1.036'; ForceText = $true },
    @{ Cell = 'E37'; Value = '  -1.51%  '; ForceText = $false },
    @{ Cell = 'D38'; Value = '0.2870'; ForceText = $true },
    @{ Cell = 'E38'; Value = '  +3.03%  '; ForceText = $false },
    @{ Cell = 'B39'; Value = 'Stellar'; ForceText = $false },
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; ForceText = $false },
    @{ Cell = 'D39'; Value = '0.09572'; ForceText = $true },
    @{ Cell = 'E39'; Value = '  +4.79%  '; ForceText = $false },
    @{ Cell = 'B40'; Value = 'FraxShare'; ForceText = $false },
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; ForceText = $false },
    @{ Cell = 'D40'; Value = '10.89'; ForceText = $true },
    @{ Cell = 'E40'; Value = '  +0.40%  '; ForceText = $false },
    @{ Cell = 'D41'; Value = '0.8234'; ForceText = $true },
    @{ Cell = 'E41'; Value = '  +2.60%  '; ForceText = $false },
    @{ Cell = 'D42'; Value = '13.99'; ForceText = $true },
    @{ Cell = 'E42'; Value = '  -1.17%  '; ForceText = $false },
    @{ Cell = 'D43'; Value = '1.481'; ForceText = $true },
    @{ Cell = 'E43'; Value = '  +0.58%  '; ForceText = $false },
    @{ Cell = 'D44'; Value = '17.26'; ForceText = $true },
    @{ Cell = 'E44'; Value = '  -1.96%  '; ForceText = $false },
    @{ Cell = 'D45'; Value = '2.679'; ForceText = $true },
    @{ Cell = 'E45'; Value = '  +0.62%  '; ForceText = $false },
    @{ Cell = 'D46'; Value = '0.7379'; ForceText = $true },
    @{ Cell = 'E46'; Value = '  +2.01%  '; ForceText = $false },
    @{ Cell = 'D47'; Value = '4.247'; ForceText = $true },
    @{ Cell = 'E47'; Value = '  -0.38%  '; ForceText = $false },
    @{ Cell = 'D48'; Value = '1.396'; ForceText = $true },
    @{ Cell = 'E48'; Value = '  -1.38%  '; ForceText = $false },
    @{ Cell = 'D49'; Value = '0.08773'; ForceText = $true },
    @{ Cell = 'E49'; Value = '  +8.59%  '; ForceText = $false },
    @{ Cell = 'D50'; Value = '1.001'; ForceText = $true },
    @{ Cell = 'E50'; Value = '  +0.17%  '; ForceText = $false },
    @{ Cell = 'D51'; Value = '139.07'; ForceText = $true },
    @{ Cell = 'E51'; Value = '  -0.06%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $cell.NumberFormat = '@'
        $cell.Value = $u.Value
        $cell.Style = 'Normal'
    } else {
        $cell.Value = $u.Value
    }
}
